$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.692.93"
$ws.Range("E2").Value = "  -1.20%  "

$ws.Range("D3").Value = "3.777.45"

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.13"
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.25"
$ws.Range("E6").Value = "  +0.00%  "

$ws.Range("D7").Value = "3.773.55"
$ws.Range("E7").Value = "  +0.31%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  -0.31%  "

$ws.Range("E10").Value = "  -0.55%  "

$ws.Range("E11").Value = "  -1.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  -2.62%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.02"
$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("D15").Value = "4.411.36"
$ws.Range("E15").Value = "  +0.39%  "

$ws.Range("D16").Value = "3.775.89"
$ws.Range("E16").Value = "  +0.49%  "

$ws.Range("D17").Value = "67.623.93"
$ws.Range("E17").Value = "  -1.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.38"
$ws.Range("E18").Value = "  +2.42%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.99"
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("E20").Value = "  -0.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("E21").Value = "  -6.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "457.31"
$ws.Range("E22").Value = "  -1.80%  "

$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("E24").Value = "  +3.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.18"
$ws.Range("E25").Value = "  -1.11%  "

$ws.Range("E26").Value = "  -0.65%  "

$ws.Range("E27").Value = "  -2.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -0.10%  "

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("E30").Value = "  -0.41%  "

$ws.Range("E31").Value = "  +3.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.23"
$ws.Range("E32").Value = "  -1.47%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.64"
$ws.Range("E33").Value = "  -1.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.12"
$ws.Range("E34").Value = "  -1.14%  "

$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("D36").Value = "3.729.26"
$ws.Range("E36").Value = "  +0.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0999"
$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("E38").Value = "  -2.11%  "

$ws.Range("E39").Value = "  -0.98%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.995"
$ws.Range("E40").Value = "  -0.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.75"
$ws.Range("E41").Value = "  -0.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.25"
$ws.Range("E44").Value = "  +2.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.29"
$ws.Range("E45").Value = "  +3.10%  "

$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.81"
$ws.Range("E47").Value = "  +3.08%  "

$ws.Range("E48").Value = "  -2.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "389.06"
$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("E50").Value = "  -5.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.54"
$ws.Range("E51").Value = "  -2.79%  "
